$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B6").Value = "/Users/nikola/pyrecodes_business/Example 5_business/water_distribution_network/rewet_results"
$ws.Range("B7").Value = "/Users/nikola/pyrecodes_business/Example 5_business/water_distribution_network/rewet_temp"
$ws.Range("B15").Value = "/Users/nikola/pyrecodes_business/Example 5_business/water_distribution_network/waterNetwork.inp"
$ws.Range("B19").Value = "/Users/nikola/pyrecodes_business/Example 5_business/water_distribution_network/rewet_temp/list.xlsx"
$ws.Range("B20").Value = "/Users/nikola/pyrecodes_business/Example 5_business/water_distribution_network/rewet_temp"
$ws.Range("B47").Value = "/Users/nikola/pyrecodes_business/env_pyrecodes/lib/python3.9/site-packages/rewet/examples/Net3/config.txt"
